# Add the new "DC Unit Loading Details" mini-table (column I, rows 1-3) to
# both worksheets, reusing the existing header/left-aligned cell styles, then
# restore the selection the author left the sheets in after reopening.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Add Devices Loop A" ---
$ws1 = $wb.Worksheets.Item("Add Devices Loop A")

# I1 : header style (same as row 7 headers, e.g. A7)
$ws1.Range("A7").Copy()
$ws1.Range("I1").PasteSpecial(-4122)
$ws1.Range("I1").Value = "DC Unit Loading Details Name"

# I2 : body style (same as row 8 cells, e.g. A8)
$ws1.Range("A8").Copy()
$ws1.Range("I2").PasteSpecial(-4122)
$ws1.Range("I2").Value = "Current (DC Units)"

# I3 : body style (same as row 8 cells, e.g. A8)
$ws1.Range("A8").Copy()
$ws1.Range("I3").PasteSpecial(-4122)
$ws1.Range("I3").Value = "Current (worst case)"

$ws1.Select()
$excel.Goto($ws1.Range("I1:I3"), $true)

# --- Sheet 2: "Other Devices Loop A" ---
$ws2 = $wb.Worksheets.Item("Other Devices Loop A")

$ws2.Range("A7").Copy()
$ws2.Range("I1").PasteSpecial(-4122)
$ws2.Range("I1").Value = "DC Unit Loading Details Name"

$ws2.Range("A8").Copy()
$ws2.Range("I2").PasteSpecial(-4122)
$ws2.Range("I2").Value = "Current (DC Units)"

$ws2.Range("A8").Copy()
$ws2.Range("I3").PasteSpecial(-4122)
$ws2.Range("I3").Value = "Current (worst case)"

# Selecting sheet2 last makes it the active/visible tab when the workbook is
# saved, matching the original workbook (activeTab=1, "Other Devices Loop A"
# tabbed as the selected sheet).
$ws2.Select()
$ws2.Range("I1:I3").Select()

$excel.CutCopyMode = $false
